$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.092.30'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.846.98'
$ws.Range('E3').Value = '  +2.46%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''360.83'
$ws.Range('E5').Value = '  +6.46%  '
$ws.Range('D6').Value = '''113.05'
$ws.Range('E6').Value = '  -1.96%  '
$ws.Range('E7').Value = '  +4.72%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '''0.605'
$ws.Range('E9').Value = '  +5.22%  '
$ws.Range('D10').Value = '''41.56'
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').Value = '''0.0865'
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('D12').Value = '''20.27'
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('E13').Value = '  +1.33%  '
$ws.Range('D14').Value = '''7.83'
$ws.Range('E14').Value = '  +3.29%  '
$ws.Range('D15').Value = '3.298.25'
$ws.Range('E15').Value = '  +2.82%  '
$ws.Range('D16').Value = '2.843.96'
$ws.Range('E16').Value = '  +2.25%  '
$ws.Range('D17').Value = '''0.933'
$ws.Range('E17').Value = '  +6.55%  '
$ws.Range('D18').Value = '52.032.97'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = '''7.63'
$ws.Range('E19').Value = '  +9.45%  '
$ws.Range('E20').Value = '  -1.10%  '
$ws.Range('D21').Value = '''13.58'
$ws.Range('E21').Value = '  +2.96%  '
$ws.Range('D22').Value = '0.0₃0999'
$ws.Range('E22').Value = '  +2.37%  '
$ws.Range('D23').Value = '''70.47'
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('D24').Value = '''269.86'
$ws.Range('E24').Value = '  -2.18%  '
$ws.Range('E25').Value = '  +5.11%  '
$ws.Range('D26').Value = '''27.19'
$ws.Range('E26').Value = '  +1.95%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D28').Value = '''10.40'
$ws.Range('E28').Value = '  +2.44%  '
$ws.Range('E29').Value = '  +1.47%  '
$ws.Range('D30').Value = '''53.46'
$ws.Range('D31').Value = '''35.72'
$ws.Range('E31').Value = '  +3.24%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '''0.140'
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('B33').Value = 'VeChain'
$ws.Range('C33').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D33').Value = '''0.0474'
$ws.Range('E33').Value = '  +25.46%  '
$ws.Range('D34').Value = '''5.96'
$ws.Range('E34').Value = '  +4.59%  '
$ws.Range('E35').Value = '  +12.38%  '
$ws.Range('D36').Value = '''0.0848'
$ws.Range('E36').Value = '  +3.65%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E38').Value = '  +3.55%  '
$ws.Range('D39').Value = '''2.07'
$ws.Range('E39').Value = '  -0.97%  '
$ws.Range('D40').Value = '''18.65'
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('E41').Value = '  +1.71%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '''23.56'
$ws.Range('E42').Value = '  +1.98%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''2.55'
$ws.Range('E43').Value = '  -4.19%  '
$ws.Range('D44').Value = '''125.25'
$ws.Range('E45').Value = '  -4.20%  '
$ws.Range('D46').Value = '''3.43'
$ws.Range('E46').Value = '  +3.90%  '
$ws.Range('D47').Value = '2.114.42'
$ws.Range('E47').Value = '  +2.31%  '
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('D49').Value = '''6.04'
$ws.Range('E49').Value = '  +9.07%  '
$ws.Range('D50').Value = '''0.984'
$ws.Range('E50').Value = '  +12.17%  '
$ws.Range('D51').Value = '''62.11'
$ws.Range('E51').Value = '  +4.88%  '
